$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at row 6 to make room for the new draws (14 July 2017 and 11 July 2017),
# shifting the existing data (and all subsequent rows) down by two rows.
$ws.Rows("6:7").Insert()

# New row for 14 July 2017
$ws.Range("A6").Value = "14 July 2017"
$ws.Range("B6").Value = 21
$ws.Range("C6").Value = 28
$ws.Range("D6").Value = 29
$ws.Range("E6").Value = 31
$ws.Range("F6").Value = 35
$ws.Range("G6").Value = 8

# New row for 11 July 2017
$ws.Range("A7").Value = "11 July 2017"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 15
$ws.Range("D7").Value = 19
$ws.Range("E7").Value = 27
$ws.Range("F7").Value = 34
$ws.Range("G7").Value = 19
